$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H17").Value = 3060.25
$ws_ALC.Range("J17").Value = 3060.25
$ws_ALC.Range("L17").Value = 9180.75
$ws_ALC.Range("N17").Value = -9516.75
$ws_ALC.Range("H40").Value = 1585.4546
$ws_ALC.Range("J40").Value = 2091.4285
$ws_ALC.Range("L40").Value = 2091.4285
$ws_ALC.Range("N40").Value = -2441.4285
$ws_ALC.Range("H62").Value = 250
$ws_ALC.Range("I62").Value = 250
$ws_ALC.Range("K62").Value = 250
$ws_ALC.Range("M62").Value = 374
$ws_ALC.Range("H65").Value = 250
$ws_ALC.Range("I65").Value = 250
$ws_ALC.Range("K65").Value = 1250
$ws_ALC.Range("M65").Value = 1870
$ws_ALC.Range("H121").Value = 1938.0769
$ws_ALC.Range("J121").Value = 1938.0769
$ws_ALC.Range("L121").Value = 5814.2307
$ws_ALC.Range("N121").Value = -9308.2307
$ws_ALC.Range("H132").Value = 3130.5186
$ws_ALC.Range("I132").Value = 1654.2
$ws_ALC.Range("J132").Value = 3998.9412
$ws_ALC.Range("K132").Value = 4962.6
$ws_ALC.Range("L132").Value = 11996.8236
$ws_ALC.Range("M132").Value = -2432.6
$ws_ALC.Range("N132").Value = -17056.8236
$ws_ALC.Range("H137").Value = 2563
$ws_ALC.Range("I137").Value = 1185.8572
$ws_ALC.Range("K137").Value = 3557.5716
$ws_ALC.Range("M137").Value = -1007.5716
$ws_ALC.Range("H138").Value = 2411.3208
$ws_ALC.Range("I138").Value = 1766.8889
$ws_ALC.Range("J138").Value = 2543.1365
$ws_ALC.Range("K138").Value = 5300.6667
$ws_ALC.Range("L138").Value = 7629.4095
$ws_ALC.Range("M138").Value = -160.6666999999998
$ws_ALC.Range("N138").Value = -17909.4095
$ws_ALC.Range("H141").Value = 5152.8887
$ws_ALC.Range("I141").Value = 5172
$ws_ALC.Range("K141").Value = 15516
$ws_ALC.Range("M141").Value = -10336

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 4261.2856
$ws_ARM.Range("I32").Value = 4269
$ws_ARM.Range("K32").Value = 4269
$ws_ARM.Range("M32").Value = -3982
$ws_ARM.Range("H35").Value = 7800
$ws_ARM.Range("I35").Value = 0
$ws_ARM.Range("K35").Value = 0
$ws_ARM.Range("M35").ClearContents()
$ws_ARM.Range("H61").Value = 1614.4667
$ws_ARM.Range("I61").Value = 972.5
$ws_ARM.Range("J61").Value = 2898.4
$ws_ARM.Range("K61").Value = 972.5
$ws_ARM.Range("L61").Value = 2898.4
$ws_ARM.Range("M61").Value = -760.5
$ws_ARM.Range("N61").Value = -3322.4
$ws_ARM.Range("H122").Value = 1179
$ws_ARM.Range("I122").Value = 1122.6364
$ws_ARM.Range("K122").Value = 3367.9092
$ws_ARM.Range("M122").Value = -917.9092000000001
$ws_ARM.Range("H132").Value = 2698.3447
$ws_ARM.Range("I132").Value = 1960.8948
$ws_ARM.Range("K132").Value = 5882.6844
$ws_ARM.Range("M132").Value = -3352.6844
$ws_ARM.Range("H136").Value = 1614.4667
$ws_ARM.Range("I136").Value = 972.5
$ws_ARM.Range("J136").Value = 2898.4
$ws_ARM.Range("K136").Value = 2917.5
$ws_ARM.Range("L136").Value = 8695.200000000001
$ws_ARM.Range("M136").Value = -367.5
$ws_ARM.Range("N136").Value = -13795.2

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H56").Value = 30000
$ws_BSM.Range("J56").Value = 30000
$ws_BSM.Range("L56").Value = 30000
$ws_BSM.Range("N56").Value = -31478
$ws_BSM.Range("H99").Value = 1995
$ws_BSM.Range("I99").Value = 1995
$ws_BSM.Range("K99").Value = 1995
$ws_BSM.Range("M99").Value = -497
$ws_BSM.Range("H105").Value = 6875
$ws_BSM.Range("I105").Value = 6875
$ws_BSM.Range("K105").Value = 6875
$ws_BSM.Range("M105").Value = -5128

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 1806.8823
$ws_CRP.Range("I31").Value = 2029.0834
$ws_CRP.Range("K31").Value = 2029.0834
$ws_CRP.Range("M31").Value = -1734.0834
$ws_CRP.Range("H34").Value = 1806.8823
$ws_CRP.Range("I34").Value = 2029.0834
$ws_CRP.Range("K34").Value = 2029.0834
$ws_CRP.Range("M34").Value = -1827.0834
$ws_CRP.Range("H99").Value = 3894.6667
$ws_CRP.Range("I99").Value = 4113.8
$ws_CRP.Range("J99").Value = 2799
$ws_CRP.Range("K99").Value = 4113.8
$ws_CRP.Range("L99").Value = 2799
$ws_CRP.Range("M99").Value = -2615.8
$ws_CRP.Range("N99").Value = -5795
$ws_CRP.Range("H126").Value = 3894.6667
$ws_CRP.Range("I126").Value = 4113.8
$ws_CRP.Range("J126").Value = 2799
$ws_CRP.Range("K126").Value = 12341.4
$ws_CRP.Range("L126").Value = 8397
$ws_CRP.Range("M126").Value = -9871.400000000001
$ws_CRP.Range("N126").Value = -13337
$ws_CRP.Range("H132").Value = 2762.8
$ws_CRP.Range("I132").Value = 2244.7
$ws_CRP.Range("J132").Value = 3799
$ws_CRP.Range("K132").Value = 6734.099999999999
$ws_CRP.Range("L132").Value = 11397
$ws_CRP.Range("M132").Value = -4204.099999999999
$ws_CRP.Range("N132").Value = -16457

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H102").Value = 4648.143
$ws_GSM.Range("I102").Value = 3304.6
$ws_GSM.Range("K102").Value = 3304.6
$ws_GSM.Range("M102").Value = -1682.6
$ws_GSM.Range("H113").Value = 3159.6
$ws_GSM.Range("I113").Value = 2932.6667
$ws_GSM.Range("K113").Value = 2932.6667
$ws_GSM.Range("M113").Value = -762.6667000000002
$ws_GSM.Range("H122").Value = 1342.1111
$ws_GSM.Range("I122").Value = 1342.1111
$ws_GSM.Range("J122").Value = 0
$ws_GSM.Range("K122").Value = 4026.3333
$ws_GSM.Range("L122").Value = 0
$ws_GSM.Range("M122").Value = -1576.3333
$ws_GSM.Range("N122").ClearContents()
$ws_GSM.Range("H126").Value = 1148.8334
$ws_GSM.Range("I126").Value = 978.8
$ws_GSM.Range("K126").Value = 2936.4
$ws_GSM.Range("M126").Value = -466.3999999999996
$ws_GSM.Range("H132").Value = 2061.5293
$ws_GSM.Range("I132").Value = 1311.8
$ws_GSM.Range("J132").Value = 2373.9167
$ws_GSM.Range("K132").Value = 3935.4
$ws_GSM.Range("L132").Value = 7121.750100000001
$ws_GSM.Range("M132").Value = -1405.4
$ws_GSM.Range("N132").Value = -12181.7501

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H68").Value = 0
$ws_LTW.Range("I68").Value = 0
$ws_LTW.Range("K68").Value = 0
$ws_LTW.Range("M68").ClearContents()
$ws_LTW.Range("H71").Value = 0
$ws_LTW.Range("I71").Value = 0
$ws_LTW.Range("K71").Value = 0
$ws_LTW.Range("M71").ClearContents()
$ws_LTW.Range("H122").Value = 3614.5557
$ws_LTW.Range("I122").Value = 3337.4666
$ws_LTW.Range("J122").Value = 5000
$ws_LTW.Range("K122").Value = 10012.3998
$ws_LTW.Range("L122").Value = 15000
$ws_LTW.Range("M122").Value = -7562.399800000001
$ws_LTW.Range("N122").Value = -19900
$ws_LTW.Range("H136").Value = 3181
$ws_LTW.Range("I136").Value = 2749.625
$ws_LTW.Range("J136").Value = 4331.3335
$ws_LTW.Range("K136").Value = 8248.875
$ws_LTW.Range("L136").Value = 12994.0005
$ws_LTW.Range("M136").Value = -5698.875
$ws_LTW.Range("N136").Value = -18094.0005

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H107").Value = 666.3333
$ws_WVR.Range("I107").Value = 666.3333
$ws_WVR.Range("K107").Value = 1998.9999
$ws_WVR.Range("M107").Value = -78.99990000000003
$ws_WVR.Range("H132").Value = 2827.739
$ws_WVR.Range("J132").Value = 3637.2307
$ws_WVR.Range("L132").Value = 10911.6921
$ws_WVR.Range("N132").Value = -15971.6921
$ws_WVR.Range("H136").Value = 2113.8823
$ws_WVR.Range("J136").Value = 2590.5557
$ws_WVR.Range("L136").Value = 7771.6671
$ws_WVR.Range("N136").Value = -12871.6671
